# Refresh the cryptos list (prices in column D, 1h change in column E)
# as published by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (e.g. '37.180.40' uses '.' for both
# thousands and decimals). For entries that would otherwise parse as a plain
# number, prefix with an apostrophe so Excel keeps storing them as text -
# same as the original cells - then restore the 'Normal' style so the
# leading apostrophe doesn't leave a quote-prefix format behind.

$ws.Range('D2').Value = '37.180.40'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '2.004.63'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'" + '258.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.70%  '
$ws.Range('D6').Value = "'" + '0.611'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.72%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'" + '56.48'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.07%  '
$ws.Range('E9').Value = '  -3.34%  '
$ws.Range('E10').Value = '  -4.88%  '
$ws.Range('E11').Value = '  -3.09%  '
$ws.Range('D12').Value = '2.300.36'
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('E13').Value = '  -6.28%  '
$ws.Range('D14').Value = "'" + '21.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('E15').Value = '  -7.21%  '
$ws.Range('E16').Value = '  -5.18%  '
$ws.Range('D17').Value = '2.015.59'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '37.236.48'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = "'" + '70.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').Value = '0.0₃0833'
$ws.Range('E20').Value = '  -3.54%  '
$ws.Range('D21').Value = "'" + '233.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = "'" + '2.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').Value = "'" + '165.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').Value = "'" + '8.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.99%  '
$ws.Range('D28').Value = "'" + '19.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('E29').Value = '  -7.46%  '
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('E32').Value = '  -4.57%  '
$ws.Range('E33').Value = '  -5.10%  '
$ws.Range('D34').Value = "'" + '4.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('E35').Value = '  -5.51%  '
$ws.Range('D36').Value = "'" + '3.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.46%  '
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = "'" + '5.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.438.03'
$ws.Range('E42').Value = '  +4.25%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = "'" + '0.0210'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = "'" + '0.0923'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.90%  '
$ws.Range('D45').Value = "'" + '89.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('D46').Value = "'" + '15.59'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.77%  '
$ws.Range('E47').Value = '  -3.57%  '
$ws.Range('D48').Value = "'" + '2.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('D49').Value = "'" + '6.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.30%  '
$ws.Range('D50').Value = '2.192.24'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('E51').Value = '  -10.06%  '
